$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 922, shifting the existing data (old rows 922:979)
# down to rows 927:984. This corresponds to a new week of price data being
# prepended to this block of "Pera" records.
$ws.Rows("922:926").Insert()

$ws.Range("A922").Value2 = 6
$ws.Range("B922").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C922").Value2 = "Metropolitana"
$ws.Range("D922").Value2 = 44585
$ws.Range("E922").Value2 = 13
$ws.Range("F922").Value2 = "Fruta"
$ws.Range("G922").Value2 = 100104
$ws.Range("H922").Value2 = "Frutos de pepita"
$ws.Range("I922").Value2 = 100104005
$ws.Range("J922").Value2 = "Pera"
$ws.Range("K922").Value2 = "Bartlett de verano"
$ws.Range("L922").Value2 = "Primera"
$ws.Range("M922").Value2 = 20
$ws.Range("N922").Value2 = 180000
$ws.Range("O922").Value2 = 180000
$ws.Range("P922").Value2 = 180000
$ws.Range("Q922").Value2 = "`$/bins (450 kilos)"
$ws.Range("R922").Value2 = "Región de O'Higgins"
$ws.Range("S922").Value2 = 400
$ws.Range("T922").Value2 = 450
$ws.Range("A923").Value2 = 6
$ws.Range("B923").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C923").Value2 = "Metropolitana"
$ws.Range("D923").Value2 = 44585
$ws.Range("E923").Value2 = 13
$ws.Range("F923").Value2 = "Fruta"
$ws.Range("G923").Value2 = 100104
$ws.Range("H923").Value2 = "Frutos de pepita"
$ws.Range("I923").Value2 = 100104005
$ws.Range("J923").Value2 = "Pera"
$ws.Range("K923").Value2 = "Bartlett de verano"
$ws.Range("L923").Value2 = "Segunda"
$ws.Range("M923").Value2 = 18
$ws.Range("N923").Value2 = 140000
$ws.Range("O923").Value2 = 140000
$ws.Range("P923").Value2 = 140000
$ws.Range("Q923").Value2 = "`$/bins (450 kilos)"
$ws.Range("R923").Value2 = "Región de O'Higgins"
$ws.Range("S923").Value2 = 311
$ws.Range("T923").Value2 = 450
$ws.Range("A924").Value2 = 6
$ws.Range("B924").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C924").Value2 = "Metropolitana"
$ws.Range("D924").Value2 = 44585
$ws.Range("E924").Value2 = 13
$ws.Range("F924").Value2 = "Fruta"
$ws.Range("G924").Value2 = 100104
$ws.Range("H924").Value2 = "Frutos de pepita"
$ws.Range("I924").Value2 = 100104005
$ws.Range("J924").Value2 = "Pera"
$ws.Range("K924").Value2 = "Favorita De Clapp"
$ws.Range("L924").Value2 = "Primera"
$ws.Range("M924").Value2 = 20
$ws.Range("N924").Value2 = 150000
$ws.Range("O924").Value2 = 150000
$ws.Range("P924").Value2 = 150000
$ws.Range("Q924").Value2 = "`$/bins (450 kilos)"
$ws.Range("R924").Value2 = "Región de O'Higgins"
$ws.Range("S924").Value2 = 333
$ws.Range("T924").Value2 = 450
$ws.Range("A925").Value2 = 6
$ws.Range("B925").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C925").Value2 = "Metropolitana"
$ws.Range("D925").Value2 = 44585
$ws.Range("E925").Value2 = 13
$ws.Range("F925").Value2 = "Fruta"
$ws.Range("G925").Value2 = 100104
$ws.Range("H925").Value2 = "Frutos de pepita"
$ws.Range("I925").Value2 = 100104005
$ws.Range("J925").Value2 = "Pera"
$ws.Range("K925").Value2 = "Favorita De Clapp"
$ws.Range("L925").Value2 = "Segunda"
$ws.Range("M925").Value2 = 20
$ws.Range("N925").Value2 = 100000
$ws.Range("O925").Value2 = 100000
$ws.Range("P925").Value2 = 100000
$ws.Range("Q925").Value2 = "`$/bins (450 kilos)"
$ws.Range("R925").Value2 = "Región de O'Higgins"
$ws.Range("S925").Value2 = 222
$ws.Range("T925").Value2 = 450
$ws.Range("A926").Value2 = 6
$ws.Range("B926").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C926").Value2 = "Metropolitana"
$ws.Range("D926").Value2 = 44585
$ws.Range("E926").Value2 = 13
$ws.Range("F926").Value2 = "Fruta"
$ws.Range("G926").Value2 = 100104
$ws.Range("H926").Value2 = "Frutos de pepita"
$ws.Range("I926").Value2 = 100104005
$ws.Range("J926").Value2 = "Pera"
$ws.Range("K926").Value2 = "Salvador Izquierdo"
$ws.Range("L926").Value2 = "Primera"
$ws.Range("M926").Value2 = 18
$ws.Range("N926").Value2 = 130000
$ws.Range("O926").Value2 = 130000
$ws.Range("P926").Value2 = 130000
$ws.Range("Q926").Value2 = "`$/bins (450 kilos)"
$ws.Range("R926").Value2 = "Región de O'Higgins"
$ws.Range("S926").Value2 = 289
$ws.Range("T926").Value2 = 450
